# Generate Report for Handback
# - Flip the "Ready for handoff" status to "Handed back: in sync with en-US"
#   everywhere it appears (Overview!E/F, zh-cn!C, de-de!C).
# - Stamp the "Latest Handback DateTime" for both target-language sheets.
# - Populate "Latest Target File" / "Latest Handback File" for each row,
#   with a hyperlink on the target-file cell (mirrors the existing source
#   hyperlink in column A).
# - Widen the columns that now hold the longer handback status / filenames.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- 1. Status text: update every cell that currently shows the old text.
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($overview.Range("E3").Value2 -eq $oldStatus) { $overview.Range("E3").Value = $newStatus }
if ($overview.Range("F3").Value2 -eq $oldStatus) { $overview.Range("F3").Value = $newStatus }

if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($zhcn.Range("C3").Value2 -eq $oldStatus) { $zhcn.Range("C3").Value = $newStatus }

if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }
if ($dede.Range("C3").Value2 -eq $oldStatus) { $dede.Range("C3").Value = $newStatus }

# ---- 2. zh-cn: fill in target/handback file columns + handback datetime.
$zhcn.Range("I2").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$zhcn.Range("J2").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-12 06:58:47"

$zhcn.Range("I3").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$zhcn.Range("J3").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-12 06:58:47"

# ---- 3. de-de: same, plus its handback datetime differs per the HO xliff.
$dede.Range("I2").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.md"
$dede.Range("J2").Value = "37e5cb7e-861c-40ec-816c-c1383e08f148.2f84ac3df99b3a1ac1251b7f21d4be11c2849a17.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 06:58:56"

$dede.Range("I3").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.md"
$dede.Range("J3").Value = "3f333ecc-78b8-442d-9710-3b4ca4700805.ac5f4a0e070ad9215f6e7102130639df11298b2f.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 06:58:56"

# ---- 4. Hyperlink the new "Latest Target File" cells to the same source
#         doc the row's source-file hyperlink (column A) already points at.
#         Rebuild each sheet's hyperlink collection in row order so the two
#         new links land right after the column-A links they mirror.
$zhcnLink2 = $zhcn.Hyperlinks.Item(1).Address
$zhcnLink3 = $zhcn.Hyperlinks.Item(2).Address
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $zhcnLink2, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $zhcnLink2, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $zhcnLink3, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $zhcnLink3, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")

$dedeLink2 = $dede.Hyperlinks.Item(1).Address
$dedeLink3 = $dede.Hyperlinks.Item(2).Address
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $dedeLink2, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$dede.Hyperlinks.Add($dede.Range("I2"), $dedeLink2, "", "", "37e5cb7e-861c-40ec-816c-c1383e08f148.md")
$dede.Hyperlinks.Add($dede.Range("A3"), $dedeLink3, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $dedeLink3, "", "", "3f333ecc-78b8-442d-9710-3b4ca4700805.md")

# ---- 5. Widen columns to fit the new, longer values.
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40
